$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = 33209.08225364993;  D = 1901.985547483333; E = 5694.2169844 }
    3  = @{ B = 31216.5484485166;   D = 1758.351291916667; E = 5239.709536516666 }
    4  = @{ B = 33173.69380958327;  D = 1909.6054191;      E = 5701.878311616667 }
    5  = @{ B = 32351.4081061166;   D = 1889.915044066667; E = 5539.472229683333 }
    6  = @{ B = 33456.34117584994;  D = 1903.622835366667; E = 5635.476059883334 }
    7  = @{ B = 32074.0560817666;   D = 1869.6174719;      E = 5756.124315266667 }
    8  = @{ B = 33022.4210715666;   D = 1899.69930645;     E = 5610.302750683333 }
    9  = @{ B = 33294.84516929993;  D = 1879.829380483333; E = 5614.832517466666 }
    10 = @{ B = 32136.99496884993;  D = 1855.6239275;      E = 5479.308287316667 }
    11 = @{ B = 33099.4474003166;   D = 1901.421045;       E = 5657.74548565 }
    12 = @{ B = 32383.45669494993;  D = 1884.783771483333; E = 5674.559654383334 }
    13 = @{ B = 32103.15360359993;  D = 1831.75313835;     E = 5761.713906583333 }
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row].B
    $ws.Range("D$row").Value = $values[$row].D
    $ws.Range("E$row").Value = $values[$row].E
}
